$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 365-366; existing rows 365-420 shift down to 367-422.
$ws.Rows("365:366").Insert()

# New row 365: Zapallo / Camote / Primera, 2021-11-22 (serial 44522)
$ws.Cells.Item(365,1).Value  = 10
$ws.Cells.Item(365,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(365,3).Value  = "La Araucanía"
$ws.Cells.Item(365,4).Value  = 44522
$ws.Cells.Item(365,5).Value  = 9
$ws.Cells.Item(365,6).Value  = 100112045
$ws.Cells.Item(365,7).Value  = "Zapallo"
$ws.Cells.Item(365,8).Value  = "Camote"
$ws.Cells.Item(365,9).Value  = "Primera"
$ws.Cells.Item(365,10).Value = 700
$ws.Cells.Item(365,11).Value = 800
$ws.Cells.Item(365,12).Value = 800
$ws.Cells.Item(365,13).Value = 800
$ws.Cells.Item(365,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(365,15).Value = "Perú"
$ws.Cells.Item(365,16).Value = 800
$ws.Cells.Item(365,17).Value = 1
$ws.Cells.Item(365,18).Value = "Hortaliza"

# New row 366: Zapallo / Paine / 1a (guarda), 2021-11-22 (serial 44522)
$ws.Cells.Item(366,1).Value  = 10
$ws.Cells.Item(366,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(366,3).Value  = "La Araucanía"
$ws.Cells.Item(366,4).Value  = 44522
$ws.Cells.Item(366,5).Value  = 9
$ws.Cells.Item(366,6).Value  = 100112045
$ws.Cells.Item(366,7).Value  = "Zapallo"
$ws.Cells.Item(366,8).Value  = "Paine"
$ws.Cells.Item(366,9).Value  = "1a (guarda)"
$ws.Cells.Item(366,10).Value = 1200
$ws.Cells.Item(366,11).Value = 250
$ws.Cells.Item(366,12).Value = 300
$ws.Cells.Item(366,13).Value = 279
$ws.Cells.Item(366,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(366,15).Value = "Región del Maule"
$ws.Cells.Item(366,16).Value = 279
$ws.Cells.Item(366,17).Value = 1
$ws.Cells.Item(366,18).Value = "Hortaliza"
